$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number to Excel's type inference
# are pre-formatted as Text so they round-trip as literal strings, matching
# the original inline-string cell contents (e.g. "237.10" rather than 237.1).

$ws.Range("D2").Value = "97.049.66"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "3.693.08"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "659.03"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.423"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.07"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("B10").Value = "USDC"
$ws.Range("C10").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.03%  "
$ws.Range("D11").Value = "3.692.57"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("E13").Value = "  +2.12%  "
$ws.Range("E14").Value = "  +11.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.78"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.66%  "
$ws.Range("D16").Value = "4.379.88"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("D17").Value = "96.803.82"
$ws.Range("E17").Value = "  +0.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "9.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.35%  "
$ws.Range("D19").Value = "3.661.59"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.81%  "
$ws.Range("E22").Value = "  -4.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "519.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  +2.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("E27").Value = "  +23.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "101.38"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "13.40"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.56"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.18"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "645.66"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.17%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.516"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +18.82%  "
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +9.52%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.17%  "
$ws.Range("E44").Value = "  +1.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "40.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.960"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0469"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("E51").Value = "  -1.98%  "
